# Multi style and config
# - Rename Sheet1 -> Style1, Sheet2 -> Style2
# - Relabel OPERATION_CONFIG / NO_WIP_REQ headers & style labels to mixed/lower case
# - Make Style2 the active sheet/tab, update various selections

$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$style1 = $wb.Worksheets.Item("Sheet1")
$style1.Name = "Style1"

$style2 = $wb.Worksheets.Item("Sheet2")
$style2.Name = "Style2"

$opConfig = $wb.Worksheets.Item("OPERATION_CONFIG")
$noWipReq = $wb.Worksheets.Item("NO_WIP_REQ")

# --- OPERATION_CONFIG: relabel headers + style values ---
# Order matters for shared-string table append order: style, Style 1, Style 2, next operation
$opConfig.Range("A1").Value = "style"
$opConfig.Range("B1").Value = "operation"

$opConfig.Range("A2").Value = "Style 1"

$opConfig.Range("A3").Value = "Style 2"
$opConfig.Range("A4").Value = "Style 2"
$opConfig.Range("A5").Value = "Style 2"
$opConfig.Range("A6").Value = "Style 2"

$opConfig.Range("C1").Value = "next operation"

# --- NO_WIP_REQ: relabel headers + style values ---
$noWipReq.Range("A1").Value = "style"
$noWipReq.Range("B1").Value = "operation"

$noWipReq.Range("A2").Value = "Style 1"
$noWipReq.Range("A3").Value = "Style 1"

$noWipReq.Range("A4").Value = "Style 2"
$noWipReq.Range("A5").Value = "Style 2"
$noWipReq.Range("A6").Value = "Style 2"

# --- Update selections on each sheet ---
$style1.Range("E3").Select()
$style2.Range("H10").Select()
$opConfig.Range("C14").Select()
$noWipReq.Range("C5").Select()

# --- Make Style2 the active sheet/tab (matches activeTab=1) ---
$style2.Activate()
